# Fix an import issue with Excel files.
# Each row on Sheet1 should have its own unique "step N" / "result N" values
# instead of all rows sharing the same "step" / "result" text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = "step 1"
$ws.Range("B2").Value = "step 2"
$ws.Range("B3").Value = "step 3"

$ws.Range("C1").Value = "result 1"
$ws.Range("C2").Value = "result 2"
$ws.Range("C3").Value = "result 3"

# Move the active selection to E4, matching the saved view state.
$ws.Activate()
$ws.Range("E4").Select()
